$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price/Volume columns for the data rows so that
# numeric-looking strings (e.g. "1.00", "315.49") are preserved as text
# instead of being auto-converted to numbers, then restore the original
# (unstyled) cell style so formatting is unchanged.
$dataRange = $ws.Range("D2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.977.05'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '2.579.57'
$ws.Range("E3").Value = '  +2.49%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '315.49'
$ws.Range("D6").Value = '100.55'
$ws.Range("E6").Value = '  +5.35%  '
$ws.Range("D7").Value = '0.574'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").Value = '36.26'
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").Value = '0.0813'
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").Value = '7.57'
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '2.978.23'
$ws.Range("E13").Value = '  +2.61%  '
$ws.Range("E14").Value = '  -0.43%  '
$ws.Range("D15").Value = '15.74'
$ws.Range("E15").Value = '  +3.39%  '
$ws.Range("D16").Value = '2.538.02'
$ws.Range("E16").Value = '  +1.47%  '
$ws.Range("D17").Value = '0.845'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '43.085.83'
$ws.Range("E18").Value = '  +0.82%  '
$ws.Range("D19").Value = '6.87'
$ws.Range("E19").Value = '  +2.13%  '
$ws.Range("D20").Value = '12.65'
$ws.Range("E20").Value = '  -1.13%  '
$ws.Range("D21").Value = '0.0₃0970'
$ws.Range("E21").Value = '  +1.21%  '
$ws.Range("D22").Value = '69.41'
$ws.Range("E22").Value = '  -0.21%  '
$ws.Range("D23").Value = '250.13'
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("D26").Value = '27.12'
$ws.Range("E26").Value = '  +2.10%  '
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("E28").Value = '  -0.92%  '
$ws.Range("D29").Value = '40.65'
$ws.Range("E29").Value = '  -1.35%  '
$ws.Range("D30").Value = '10.32'
$ws.Range("E30").Value = '  -0.15%  '
$ws.Range("D31").Value = '5.84'
$ws.Range("E31").Value = '  -1.87%  '
$ws.Range("D32").Value = '157.66'
$ws.Range("E32").Value = '  -0.30%  '
$ws.Range("D33").Value = '3.43'
$ws.Range("E33").Value = '  +4.52%  '
$ws.Range("D34").Value = '2.13'
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("D35").Value = '0.0805'
$ws.Range("E35").Value = '  +3.49%  '
$ws.Range("D36").Value = '2.67'
$ws.Range("E36").Value = '  -0.40%  '
$ws.Range("D37").Value = '18.90'
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("E38").Value = '  +10.06%  '
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("B40").Value = 'Stellar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D40").Value = '0.119'
$ws.Range("E40").Value = '  +0.46%  '
$ws.Range("B41").Value = 'EnergySwap'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D41").Value = '23.94'
$ws.Range("E41").Value = '  +3.10%  '
$ws.Range("D42").Value = '4.07'
$ws.Range("E42").Value = '  +7.55%  '
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("E44").Value = '  +0.06%  '
$ws.Range("D45").Value = '3.27'
$ws.Range("E45").Value = '  -1.42%  '
$ws.Range("D46").Value = '2.004.21'
$ws.Range("E46").Value = '  -1.72%  '
$ws.Range("D47").Value = '8.92'
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").Value = '2.828.44'
$ws.Range("E48").Value = '  +2.48%  '
$ws.Range("D49").Value = '0.197'
$ws.Range("E49").Value = '  +2.45%  '
$ws.Range("D50").Value = '75.12'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("D51").Value = '81.95'
$ws.Range("E51").Value = '  -2.85%  '

# Restore original style on the data range
$dataRange.Style = $origStyle

